$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Populate new test case rows 24-33 with values ---
# Row 24
$ws.Range("B24").Value2 = "TC_HA_21"
$ws.Range("C24").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D24").Value2 = "Validate the functionality of the notificatons in `"Your account`" dropdown menu. "
$ws.Range("E24").Value2 = "User must be loged in."
$ws.Range("F24").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Clcik any notification icon.`n4. Click on the save button."
$ws.Range("G24").Value2 = "Not Applicable"
$ws.Range("H24").Value2 = "The color of the messagae icon should be changed with a low gradient."

# Row 25
$ws.Range("B25").Value2 = "TC_HA_22 "
$ws.Range("C25").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D25").Value2 = "Validate the display of he `"Privacy and Security`" section of `"Your account`" page."
$ws.Range("E25").Value2 = "User must be loged in."
$ws.Range("F25").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button."
$ws.Range("G25").Value2 = "Not Applicable"
$ws.Range("H25").Value2 = "1. 4 buttons should be mentioned`n        a. Your apps.`n        b. Recent Visits.`n        c. Logout on all devices.`n        d. Manage Your blocklist.`n2. Check box of friend suggestions."

# Row 26
$ws.Range("B26").Value2 = "TC_HA_23"
$ws.Range("C26").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D26").Value2 = "Validate the functionality of the `"Your apps `" button in  `"Privacy and Security`" section of `"Your account`" page."
$ws.Range("E26").Value2 = "User must be loged in."
$ws.Range("F26").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button.`n4. Click on the Your apps button."
$ws.Range("G26").Value2 = "Not Applicable"
$ws.Range("H26").Value2 = "User should navigate to the page of bulding own application."

# Row 27
$ws.Range("B27").Value2 = "TC_HA_24"
$ws.Range("C27").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D27").Value2 = "Validate the Display of the `"Recent Visits`" button in  `"Privacy and Security`" section of `"Your account`" page."
$ws.Range("E27").Value2 = "User must be loged in."
$ws.Range("F27").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button.`n4. Click on the recent visits buttons button."
$ws.Range("G27").Value2 = "Not Applicable"
$ws.Range("H27").Value2 = "1. Link of the contact support should be mentioned.`n2. table with (Time, Device, Estimated loction) should be mentioned.`n3. show additional technical details button should be mentioned"
$ws.Range("I27").Value2 = " "

# Row 28
$ws.Range("B28").Value2 = "TC_HA_25"
$ws.Range("C28").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D28").Value2 = "Validate the Display of the of `"Contact Support `" link in `"Recent Visits`" button in  `"Privacy and Security`" section of `"Your account`" page."
$ws.Range("E28").Value2 = "User must be loged in."
$ws.Range("F28").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button.`n4. Click on the recent visits buttons button.`n5. Click on the Contact & Support button."
$ws.Range("G28").Value2 = "Not Applicable"
$ws.Range("H28").Value2 = "User should navigate to the contact and support page."

# Row 29
$ws.Range("B29").Value2 = "TC_HA_26"
$ws.Range("C29").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D29").Value2 = "Validate the functionality of the `"Recent Visits`" button in  `"Privacy and Security`" section of `"Your account`" page login in different devices.."
$ws.Range("E29").Value2 = "User must be logged in in different devices."
$ws.Range("F29").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button.`n4.Login using multiple devices.`n5. Click on the recent visits buttons button.`n6. Click on the Contact & Support button."
$ws.Range("G29").Value2 = "Devices:`n1. I phone`n2. Android. `n3. Windows."
$ws.Range("H29").Value2 = "Time Device and Location should be mentioned in the table when logged in to the device."

# Row 30
$ws.Range("B30").Value2 = "TC_HA_27"
$ws.Range("C30").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D30").Value2 = "Validate the Display of the `"Recent Visits`" button in  `"Privacy and Security`" section of `"Your account`" page by selecting `"Show additional technical details`"."
$ws.Range("E30").Value2 = "User must be loged in."
$ws.Range("F30").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button.`n4. Click on the recent visits buttons button.`n5. Click on the Contact & Support button."
$ws.Range("G30").Value2 = "Not Applicable"
$ws.Range("H30").Value2 = "Table should be mentioned with 5 headings.`n1.Time`n2. Device.`n3. User Agent.`n4. Estimated Location.`n5. IP Address.`nAlong with their respective details."

# Row 31
$ws.Range("B31").Value2 = "TC_HA_28"
$ws.Range("C31").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D31").Value2 = "Validate the functionality of the `"Log out on all devices`" button in  `"Privacy and Security`" section of `"Your account`" page while being logged in only one device..."
$ws.Range("E31").Value2 = "User must be loged in."
$ws.Range("F31").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button.`n4. Click on the `"log out all devices`" button."
$ws.Range("G31").Value2 = "Not Applicable"
$ws.Range("H31").Value2 = "User should be logged out right after clicking the button."

# Row 32
$ws.Range("B32").Value2 = "TC_HA_29"
$ws.Range("C32").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D32").Value2 = "Validate the functionality of the `"Log out on all devices`" button in  `"Privacy and Security`" section of `"Your account`" page while being logged in more than one devices."
$ws.Range("E32").Value2 = "User must be loged in atleast in 2 devices."
$ws.Range("F32").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button.`n4. Click on the `"log out all devices`" button."
$ws.Range("G32").Value2 = "Not Applicable"
$ws.Range("H32").Value2 = "User shold be logged out from all the devices."

# Row 33
$ws.Range("B33").Value2 = "TC_HA_30"
$ws.Range("C33").Value2 = "(TS_004)`nHeader_Your account"
$ws.Range("D33").Value2 = "Validate the display of  `"manage your block list`" button in  `"Privacy and Security`" section of `"Your account`" page while being logged in more than one devices."
$ws.Range("E33").Value2 = "1. User must be loged in.`n2. no account should be blocked yet."
$ws.Range("F33").Value2 = "1 In home page click the top right corner.`n2. Select the `"Your account `" in the dropdown menu.`n3. Scroll down the page to Privacy & Security button.`n4. Click on the `"manage your block list`" button."
$ws.Range("G33").Value2 = "Not Applicable"
$ws.Range("H33").Value2 = "1.User should navigate to new tab along with some details about the blocked account `n2.Message should me displayed  about no one is blocked yet.`n3. input field for the Email.`n4. Block button."

# --- Fix up "G" column alignment/style: left+center, no wrap (matches style used for "Not Applicable" cells elsewhere) ---
$gRows = @(24,25,26,27,28,30,31,32,33)
foreach ($r in $gRows) {
    $g = $ws.Range("G$r")
    $g.HorizontalAlignment = -4131
    $g.VerticalAlignment = -4108
    $g.WrapText = $false
}

# Row 29s G cell keeps the center/center/wrap style (like the rest of column G header style)
$g29 = $ws.Range("G29")
$g29.HorizontalAlignment = -4108
$g29.VerticalAlignment = -4108
$g29.WrapText = $true

# --- "H" and "I" columns: left-aligned, top-valign, wrap text (matches other description cells) ---
$hRows = @(24,25,26,27,28,29,30,31,32,33)
foreach ($r in $hRows) {
    $h = $ws.Range("H$r")
    $h.HorizontalAlignment = -4131
    $h.VerticalAlignment = -4160
    $h.WrapText = $true
}
$i27 = $ws.Range("I27")
$i27.HorizontalAlignment = -4131
$i27.VerticalAlignment = -4160
$i27.WrapText = $true

# --- Set row heights to match the content-driven heights used in the sheet ---
$ws.Rows.Item(24).RowHeight = 75
$ws.Rows.Item(25).RowHeight = 90
$ws.Rows.Item(26).RowHeight = 90
$ws.Rows.Item(27).RowHeight = 90
$ws.Rows.Item(28).RowHeight = 105
$ws.Rows.Item(29).RowHeight = 120
$ws.Rows.Item(30).RowHeight = 120
$ws.Rows.Item(31).RowHeight = 105
$ws.Rows.Item(32).RowHeight = 105
$ws.Rows.Item(33).RowHeight = 105

# --- Update the active selection / view to reflect the edited area ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H34").Select()
